$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the hyperlink for the new row, pointing at the new commit diff URL
$target = "https://github.com/nguyentienminh07102004/product-management/commit/7bb55ddda4cb256a0125f3d907b4e273ec7230a9"
$location = "diff-4ea027d01c6e49507fa3f91d5a63f7466ff135dff36912874b4cf40bebe35889"
$display = "https://github.com/nguyentienminh07102004/product-management/commit/7bb55ddda4cb256a0125f3d907b4e273ec7230a9 - diff-4ea027d01c6e49507fa3f91d5a63f7466ff135dff36912874b4cf40bebe35889"

$ws.Hyperlinks.Add($ws.Range("B23"), $target, $location, "", $display)

# Re-apply B22's "Hyperlink" cell style onto B23 (reuses the existing style entry) and set the new text
$ws.Range("B22").Copy($ws.Range("B23"))
$ws.Range("B23").Value = "2. tiny MCE base "

# Update selection to match new active cell
$ws.Range("C23").Select() | Out-Null
